$p = $ppt.ActivePresentation
$s3 = $p.Slides.Item(3)
$shp = $s3.Shapes.AddShape(1, 100, 100, 200, 200)
Write-Host "Type:" $shp.Type
Write-Host "AutoShapeType:" $shp.AutoShapeType
$shp.Fill.ForeColor.SchemeColor = 5
Write-Host "SchemeColor readback:" $shp.Fill.ForeColor.SchemeColor
